$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the existing "Late" column (col 14 = N)
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab, with L14 selected
$ws.Activate()
$ws.Range("L14").Select()
